$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = $ws.Range("A2:E2")
$row2.NumberFormat = "@"

$ws.Range("A2").Value = "Álcool Isopropílico"
$ws.Range("B2").Value = "2"
$ws.Range("C2").Value = "https://www.amazon.com.br/Álcool-Isopropílico-Limpeza-Placas-Circuitos/dp/B0DGMJ9633/ref=asc_df_B0DGMJ9633/?tag=googleshopp00-20&linkCode=df0&hvadid=709986098834&hvpos=&hvnetw=g&hvrand=5511000974376961432&hvpone=&hvptwo=&hvqmt=&hvdev=c&hvdvcmdl=&hvlocint=&hvlocphy=9197371&hvtargid=pla-2364142716389&psc=1&mcid=981e11534b6c31b09f06ed5bc85e56a3&gad_source=1"
$ws.Range("D2").Value = "2024-09-30"
$ws.Range("E2").Value = "outros"

$row2.ClearFormats()
